$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the betting-data content of rows 150 (id=148) and 151 (id=149).
#    Column A (the running id) stays put; every other column (B..AC) from
#    the two rows is exchanged.
# ---------------------------------------------------------------------------

function Get-RowValues($row) {
    $cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    $vals = @{}
    foreach ($col in $cols) {
        $vals[$col] = $ws.Range($col + $row).Value2
    }
    return $vals
}

$row150 = Get-RowValues 150
$row151 = Get-RowValues 151

function Set-RowValues($row, $vals) {
    foreach ($col in $vals.Keys) {
        $ws.Range($col + $row).Value2 = $vals[$col]
    }
}

Set-RowValues 150 $row151
Set-RowValues 151 $row150

# ---------------------------------------------------------------------------
# 2) Append a new row (162) for the still-unplayed fixture id=160.
# ---------------------------------------------------------------------------

$ws.Cells.Item(162, 1).Value2 = 160
$ws.Cells.Item(162, 2).Value2 = 6803258
$ws.Cells.Item(162, 3).Value = "Czech Republic 2 Liga"
$ws.Cells.Item(162, 4).Value = "Czech Republic 2 Liga"
$ws.Cells.Item(162, 5).Value2 = 45387.54166666666
$ws.Cells.Item(162, 6).Value = "FC Sellier  Bellot Vlasim"
$ws.Cells.Item(162, 7).Value = "Viktoria Zizkov"
# H (FTHG), I (FTAG) and J (FTR) are intentionally left blank: match not played yet.
$ws.Cells.Item(162, 11).Value2 = 2
$ws.Cells.Item(162, 12).Value2 = 3.5
$ws.Cells.Item(162, 13).Value2 = 3.25
$ws.Cells.Item(162, 14).Value2 = 2.2
$ws.Cells.Item(162, 15).Value2 = 3.4
$ws.Cells.Item(162, 16).Value2 = 3
$ws.Cells.Item(162, 17).Value2 = -0.25
$ws.Cells.Item(162, 18).Value2 = 1.95
$ws.Cells.Item(162, 19).Value2 = 1.85
$ws.Cells.Item(162, 20).Value2 = 2.5
$ws.Cells.Item(162, 21).Value2 = 1.85
$ws.Cells.Item(162, 22).Value2 = 1.95
$ws.Cells.Item(162, 23).Value2 = 0
$ws.Cells.Item(162, 24).Value2 = 0
$ws.Cells.Item(162, 25).Value2 = 0
$ws.Cells.Item(162, 26).Value2 = 0
$ws.Cells.Item(162, 27).Value2 = 0
# AB (PL_AhOver) and AC (PL_AhUnder) are intentionally left blank.

# Match the look of the other rows: bold/bordered id cell and the
# custom date format for the match date cell.
$ws.Range("A161").Copy() | Out-Null
$ws.Range("A162").PasteSpecial(-4122) | Out-Null

$ws.Range("E161").Copy() | Out-Null
$ws.Range("E162").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
